$p = $ppt.ActivePresentation

# --- Slide 14: merge the three runs of the "Just as with ordinary functions..." paragraph
#     into a single run (keeping the formatting of the first run).
$s14 = $p.Slides.Item(14)
$shp14 = $s14.Shapes.Item(2)
$tr14 = $shp14.TextFrame.TextRange
$found14 = $tr14.Find("Just as with ordinary functions, you can define local variables inside ")
if ($found14 -ne $null) {
    $mergedText14 = "Just as with ordinary functions, you can define local variables inside the body of a lambda expression"
    $full14 = $tr14.Characters($found14.Start, $mergedText14.Length)
    $full14.Text = $mergedText14
}

# --- Slide 15: fix missing closing paren in the lambda parameter list.
$s15 = $p.Slides.Item(15)
$shp15 = $s15.Shapes.Item(2)
$tr15 = $shp15.TextFrame.TextRange
$found15 = $tr15.Find("((Integer x -> {")
if ($found15 -ne $null) {
    $found15.Text = "((Integer x) -> {"
}

# --- Slide 18: fix "Arrrays.asList" typo -> "Arrays.asList"
$s18 = $p.Slides.Item(18)
$shp18 = $s18.Shapes.Item(2)
$tr18 = $shp18.TextFrame.TextRange
$found18 = $tr18.Find("Arrrays.asList")
if ($found18 -ne $null) {
    $found18.Text = "Arrays.asList"
}
